# Scheduled-runner style refresh of the Leve price/profit columns (H:N) across
# all eight crafting-job sheets, reflecting newly pulled Universalis market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 679.125
$ws.Range("I6").Value = 769.7143
$ws.Range("K6").Value = 2309.1429
$ws.Range("M6").Value = -2197.1429

$ws.Range("H9").Value = 83.85714
$ws.Range("I9").Value = 51.285713
$ws.Range("K9").Value = 51.285713
$ws.Range("M9").Value = 117.714287

$ws.Range("H12").Value = 411.4
$ws.Range("I12").Value = 499.25
$ws.Range("K12").Value = 499.25
$ws.Range("M12").Value = -329.25

$ws.Range("H38").Value = 2458.2144
$ws.Range("I38").Value = 1041.5
$ws.Range("K38").Value = 3124.5
$ws.Range("M38").Value = -2752.5

$ws.Range("H43").Value = 4998
$ws.Range("I43").Value = 4998
$ws.Range("K43").Value = 4998
$ws.Range("M43").Value = -4929

$ws.Range("H58").Value = 2361.3333
$ws.Range("J58").Value = 3349.5
$ws.Range("L58").Value = 10048.5
$ws.Range("N58").Value = -10348.5

$ws.Range("H61").Value = 9000.333000000001
$ws.Range("I61").Value = 9000.5
$ws.Range("J61").Value = 9000
$ws.Range("K61").Value = 27001.5
$ws.Range("L61").Value = 27000
$ws.Range("M61").Value = -26829.5
$ws.Range("N61").Value = -27344

$ws.Range("H112").Value = 1863.8889
$ws.Range("J112").Value = 2029.0435
$ws.Range("L112").Value = 6087.1305
$ws.Range("N112").Value = -8303.130499999999

$ws.Range("H137").Value = 6312.6665
$ws.Range("I137").Value = 1756.9286
$ws.Range("J137").Value = 15424.143
$ws.Range("K137").Value = 5270.7858
$ws.Range("L137").Value = 46272.429
$ws.Range("M137").Value = -2720.7858
$ws.Range("N137").Value = -51372.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1412678.8
$ws.Range("I32").Value = 700652.5
$ws.Range("K32").Value = 700652.5
$ws.Range("M32").Value = -700365.5

$ws.Range("H105").Value = 49999.5
$ws.Range("J105").Value = 49999.5
$ws.Range("L105").Value = 49999.5
$ws.Range("N105").Value = -56987.5

$ws.Range("H110").Value = 3534.2354
$ws.Range("I110").Value = 3541.5715
$ws.Range("K110").Value = 3541.5715
$ws.Range("M110").Value = -1496.5715

$ws.Range("H132").Value = 1864.3125
$ws.Range("I132").Value = 1738.341
$ws.Range("K132").Value = 5215.022999999999
$ws.Range("M132").Value = -2685.022999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 98649.664
$ws.Range("J9").Value = 98649.664
$ws.Range("L9").Value = 98649.664
$ws.Range("N9").Value = -98985.664

$ws.Range("H99").Value = 205961.8
$ws.Range("I99").Value = 337436.34
$ws.Range("J99").Value = 8750
$ws.Range("K99").Value = 337436.34
$ws.Range("L99").Value = 8750
$ws.Range("M99").Value = -335938.34
$ws.Range("N99").Value = -11746

$ws.Range("H105").Value = 13001484
$ws.Range("I105").Value = 626011.0600000001
$ws.Range("K105").Value = 626011.0600000001
$ws.Range("M105").Value = -624264.0600000001

$ws.Range("H134").Value = 3272
$ws.Range("J134").Value = 3285.4285
$ws.Range("L134").Value = 9856.2855
$ws.Range("N134").Value = -14926.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5685327
$ws.Range("I31").Value = 2649.5
$ws.Range("J31").Value = 6948144.5
$ws.Range("K31").Value = 2649.5
$ws.Range("L31").Value = 6948144.5
$ws.Range("M31").Value = -2354.5
$ws.Range("N31").Value = -6948734.5

$ws.Range("H34").Value = 5685327
$ws.Range("I34").Value = 2649.5
$ws.Range("J34").Value = 6948144.5
$ws.Range("K34").Value = 2649.5
$ws.Range("L34").Value = 6948144.5
$ws.Range("M34").Value = -2447.5
$ws.Range("N34").Value = -6948548.5

$ws.Range("H132").Value = 1887.2188
$ws.Range("I132").Value = 1526.862
$ws.Range("J132").Value = 5370.6665
$ws.Range("K132").Value = 4580.586
$ws.Range("L132").Value = 16111.9995
$ws.Range("M132").Value = -2050.586
$ws.Range("N132").Value = -21171.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 1799.6666
$ws.Range("I54").Value = 1949.5
$ws.Range("J54").Value = 1500
$ws.Range("K54").Value = 5848.5
$ws.Range("L54").Value = 4500
$ws.Range("M54").Value = -5289.5
$ws.Range("N54").Value = -5618

$ws.Range("H57").Value = 2999
$ws.Range("I57").Value = 2999
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 8997
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -8438
$ws.Range("N57").ClearContents()

$ws.Range("H68").Value = 1889663.1
$ws.Range("I68").Value = 2124.375
$ws.Range("K68").Value = 6373.125
$ws.Range("M68").Value = -5562.125

$ws.Range("H71").Value = 1889663.1
$ws.Range("I71").Value = 2124.375
$ws.Range("K71").Value = 19119.375
$ws.Range("M71").Value = -15063.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 4978
$ws.Range("J36").Value = 2713
$ws.Range("L36").Value = 2713
$ws.Range("N36").Value = -3683

$ws.Range("H70").Value = 38467212
$ws.Range("I70").Value = 62504284
$ws.Range("K70").Value = 62504284
$ws.Range("M70").Value = -62504014

$ws.Range("H73").Value = 38467212
$ws.Range("I73").Value = 62504284
$ws.Range("K73").Value = 62504284
$ws.Range("M73").Value = -62503348

$ws.Range("H122").Value = 76929010
$ws.Range("I122").Value = 250002620
$ws.Range("K122").Value = 750007860
$ws.Range("M122").Value = -750005410

$ws.Range("H132").Value = 2765.9412
$ws.Range("I132").Value = 2846.7727
$ws.Range("J132").Value = 2617.75
$ws.Range("K132").Value = 8540.3181
$ws.Range("L132").Value = 7853.25
$ws.Range("M132").Value = -6010.3181
$ws.Range("N132").Value = -12913.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 986.625
$ws.Range("I22").Value = 558.5
$ws.Range("J22").Value = 1414.75
$ws.Range("K22").Value = 558.5
$ws.Range("L22").Value = 1414.75
$ws.Range("M22").Value = -263.5
$ws.Range("N22").Value = -2004.75

$ws.Range("H27").Value = 986.625
$ws.Range("I27").Value = 558.5
$ws.Range("J27").Value = 1414.75
$ws.Range("K27").Value = 558.5
$ws.Range("L27").Value = 1414.75
$ws.Range("M27").Value = -451.5
$ws.Range("N27").Value = -1628.75

$ws.Range("H55").Value = 632.8
$ws.Range("I55").Value = 471.2857
$ws.Range("K55").Value = 471.2857
$ws.Range("M55").Value = -298.2857

$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996

$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984

$ws.Range("H100").Value = 1181.5454
$ws.Range("I100").Value = 1111.1111
$ws.Range("J100").Value = 1498.5
$ws.Range("K100").Value = 1111.1111
$ws.Range("L100").Value = 1498.5
$ws.Range("M100").Value = -570.1111000000001
$ws.Range("N100").Value = -2580.5

$ws.Range("H136").Value = 5620.684
$ws.Range("I136").Value = 3618.7273
$ws.Range("K136").Value = 10856.1819
$ws.Range("M136").Value = -8306.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 1168888
$ws.Range("J64").Value = 1168888
$ws.Range("L64").Value = 1168888
$ws.Range("N64").Value = -1169384

$ws.Range("H67").Value = 1168888
$ws.Range("J67").Value = 1168888
$ws.Range("L67").Value = 1168888
$ws.Range("N67").Value = -1170604

$ws.Range("H113").Value = 883.5263
$ws.Range("I113").Value = 1040.8
$ws.Range("J113").Value = 293.75
$ws.Range("K113").Value = 3122.4
$ws.Range("L113").Value = 881.25
$ws.Range("M113").Value = -952.3999999999996
$ws.Range("N113").Value = -5221.25

$ws.Range("H118").Value = 120999
$ws.Range("J118").Value = 120999
$ws.Range("L118").Value = 120999
$ws.Range("N118").Value = -124313

$ws.Range("H132").Value = 3040.0278
$ws.Range("I132").Value = 2971.1785
$ws.Range("J132").Value = 3281
$ws.Range("K132").Value = 8913.5355
$ws.Range("L132").Value = 9843
$ws.Range("M132").Value = -6383.5355
$ws.Range("N132").Value = -14903

$ws.Range("H136").Value = 10313.777
$ws.Range("I136").Value = 9454.5
$ws.Range("K136").Value = 28363.5
$ws.Range("M136").Value = -25813.5
